# SectorGroup.xlsx update:
# The codeforiati: columns were re-ordered -- what used to be the
# "category-name"/"group-name" pair (columns D/E) is swapped, and what
# used to be the "group-code"/"category-code" pair (columns F/G) is
# swapped too (this also flips the D1/E1 and F1/G1 header captions).
#
# Net effect for every row in the used range: swap column D <-> E and
# swap column F <-> G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()
    $gVal = $gCell.Value()

    # D <-> E : these are always non-numeric captions/labels, so a plain
    # value assignment keeps them stored as text (matches original type).
    $dCell.Value = $eVal
    $eCell.Value = $dVal

    # F <-> G : these hold numeric-looking codes ("110", "111", ...).
    # Force the Text number format before assigning so the engine keeps
    # them as shared-string text (as in the source file) instead of
    # silently re-typing them as numbers, then restore the cell style so
    # no stray formatting is left behind.
    $fCell.NumberFormat = "@"
    $fCell.Value = $gVal
    $fCell.Style = "Normal"

    $gCell.NumberFormat = "@"
    $gCell.Value = $fVal
    $gCell.Style = "Normal"
}
